$wb = $excel.ActiveWorkbook

# --- Sheets before the edit: "2010 and 2010-18" (1), "2000-09" (2) ---
$wsBaseline = $wb.Worksheets.Item(1)   # "2010 and 2010-18"
$wsOld09    = $wb.Worksheets.Item(2)   # "2000-09" (will stay "2000-09", just moves position)

# --- Insert a new blank worksheet "Sheet1" between them ---
$newSheet = $wb.Worksheets.Add($null, $wsBaseline)
$newSheet.Name = "Sheet1"

# Re-resolve the "2000-09" worksheet by name (keeps working regardless of index shuffles)
$ws09 = $wb.Worksheets.Item("2000-09")
$ws18 = $wb.Worksheets.Item("2010 and 2010-18")

# =========================================================================
# 1) "2000-09" sheet: insert a new data row (row 4), a styled blank
#    separator row (row 5, pushing the old header row from 5 -> 6), and
#    widen column R (18).
# =========================================================================
$ws09.Rows.Item(5).Insert()

$ws09.Range("A4").Value = "CW3M"
$ws09.Range("B4").Value = "Baseline_2000-09_C194"
$ws09.Range("C4").Value = "2000-09"

$ws09.Range("D4").Value = 931.63393560000009
$ws09.Range("D4").NumberFormat = "0.00"
$ws09.Range("D4").Interior.ColorIndex = 6

$ws09.Range("E4").Value = 1609.2949586000002
$ws09.Range("E4").NumberFormat = "0.00"

$ws09.Range("F4").Value = 1.0215697000000001
$ws09.Range("F4").NumberFormat = "0.00"

$ws09.Range("G4").Value = 305.74331049999995
$ws09.Range("G4").NumberFormat = "0.00"
$ws09.Range("G4").Interior.ColorIndex = 6

$ws09.Range("H4").Value = 9.3183378000000001
$ws09.Range("H4").NumberFormat = "0.00"

$ws09.Range("I4").Value = 8.1064159
$ws09.Range("I4").NumberFormat = "0.00"

$ws09.Range("J4").Value = 7.7646284999999988
$ws09.Range("J4").NumberFormat = "0.00"

$ws09.Range("K4").Value = 654.30700679999995
$ws09.Range("K4").NumberFormat = "0.00"
$ws09.Range("K4").Interior.ColorIndex = 6

$ws09.Range("L4").Value = 61.771183299999997
$ws09.Range("L4").NumberFormat = "0.00"

$ws09.Range("M4").Value = 1196.3218261999998
$ws09.Range("M4").NumberFormat = "0.00"
$ws09.Range("M4").Interior.ColorIndex = 6

$ws09.Range("N4").Value = 945.03428959999997
$ws09.Range("N4").NumberFormat = "0.00"
$ws09.Range("N4").Interior.ColorIndex = 6

$ws09.Range("O4").Value = 5406.1934815999994
$ws09.Range("O4").NumberFormat = "0"
$ws09.Range("O4").Interior.ColorIndex = 6

$ws09.Range("P4").Value = 25979.647461100001
$ws09.Range("P4").NumberFormat = "0"

$ws09.Range("Q4").Value = 0.080406299999999972
$ws09.Range("Q4").NumberFormat = "0.00"

$ws09.Range("R4").Value = -0.000016500000000000015
$ws09.Range("R4").NumberFormat = "0.000000"

$ws09.Range("S4").Value = "2000-09"

# Blank styled separator row (row 5) - numeric formats only, no values
$ws09.Range("D5:N5").NumberFormat = "0.00"
$ws09.Range("O5:P5").NumberFormat = "0"
$ws09.Range("Q5").NumberFormat = "0.00"
$ws09.Range("R5").NumberFormat = "0.000000"

# Widen column R (18)
$ws09.Columns.Item(18).ColumnWidth = 9.72

# Update the saved selection on this sheet
$ws09.Activate()
$ws09.Range("D5").Select()

# =========================================================================
# 2) "2010 and 2010-18" sheet: append a new data row (row 64).
# =========================================================================
$ws18.Range("A64").Value = "CW3M"
$ws18.Range("B64").Value = "Baseline_2010-current_2010-19_C195"
$ws18.Range("C64").Value = "2010-18"

$ws18.Range("D64").Value = 1070.3662515555557
$ws18.Range("D64").NumberFormat = "0.00"
$ws18.Range("D64").Interior.ColorIndex = 6

$ws18.Range("E64").Value = 1763.5263265555557
$ws18.Range("E64").NumberFormat = "0.00"

$ws18.Range("F64").Value = 0.999942
$ws18.Range("F64").NumberFormat = "0.00"

$ws18.Range("G64").Value = 305.6782124444444
$ws18.Range("G64").NumberFormat = "0.00"

$ws18.Range("H64").Value = 9.775355222222224
$ws18.Range("H64").NumberFormat = "0.00"

$ws18.Range("I64").Value = 6.8224234444444436
$ws18.Range("I64").NumberFormat = "0.00"

$ws18.Range("J64").Value = 8.145128999999999
$ws18.Range("J64").NumberFormat = "0.00"

$ws18.Range("K64").Value = 672.51038266666671
$ws18.Range("K64").NumberFormat = "0.00"

$ws18.Range("L64").Value = 60.018756111111117
$ws18.Range("L64").NumberFormat = "0.00"

$ws18.Range("M64").Value = 1335.0520562222218
$ws18.Range("M64").NumberFormat = "0.00"
$ws18.Range("M64").Interior.ColorIndex = 6

$ws18.Range("N64").Value = 1081.7151217777778
$ws18.Range("N64").NumberFormat = "0.00"
$ws18.Range("N64").Interior.ColorIndex = 6

$ws18.Range("O64").Value = 4576.182644333333
$ws18.Range("O64").NumberFormat = "0"

$ws18.Range("P64").Value = 27227.338324888889
$ws18.Range("P64").NumberFormat = "0"

$ws18.Range("Q64").Value = 0.27293466666666671
$ws18.Range("Q64").NumberFormat = "0.00"

$ws18.Range("R64").Value = 0.000057555555555555559
$ws18.Range("R64").NumberFormat = "0.000000"

$ws18.Range("S64").Value = "2010-18"

# Update the saved selection/view on this sheet
$ws18.Activate()
$ws18.Range("S65").Select()
